$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the "as of" snapshot time in column D for every data row (2-48) ---
$ws.Range("D2:D48").Value = 45968.473356481481

# --- Rows 19-48: newer charging-idle records replace the previous snapshot ---
$ws.Range("A19").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B19").Value = "306号直流"
$ws.Range("C19").Value = 45964.263055555559
$ws.Range("A20").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B20").Value = "004A号直流"
$ws.Range("C20").Value = 45964.528668981482
$ws.Range("A21").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B21").Value = "502号直流"
$ws.Range("C21").Value = 45965.254895833335
$ws.Range("A22").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B22").Value = "604号直流"
$ws.Range("C22").Value = 45965.565891203703
$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "402号直流"
$ws.Range("C23").Value = 45966.207546296297
$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "603号直流"
$ws.Range("C24").Value = 45966.254062499997
$ws.Range("A25").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B25").Value = "110号直流"
$ws.Range("C25").Value = 45966.540949074071
$ws.Range("A26").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B26").Value = "406号直流"
$ws.Range("C26").Value = 45966.690613425926
$ws.Range("A27").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B27").Value = "504号直流"
$ws.Range("C27").Value = 45967.035775462966
$ws.Range("A28").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B28").Value = "001A号直流"
$ws.Range("C28").Value = 45967.03800925926
$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "B03号直流"
$ws.Range("C29").Value = 45967.067662037036
$ws.Range("A30").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B30").Value = "405号直流"
$ws.Range("C30").Value = 45967.114155092589
$ws.Range("A31").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B31").Value = "404号直流"
$ws.Range("C31").Value = 45967.131412037037
$ws.Range("A32").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B32").Value = "A01号直流"
$ws.Range("C32").Value = 45967.401446759257
$ws.Range("A33").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B33").Value = "505号直流"
$ws.Range("C33").Value = 45967.507719907408
$ws.Range("A34").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B34").Value = "003B号直流"
$ws.Range("C34").Value = 45967.530300925922
$ws.Range("A35").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B35").Value = "001B号直流"
$ws.Range("C35").Value = 45967.543043981481
$ws.Range("A36").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B36").Value = "103号直流"
$ws.Range("C36").Value = 45967.554571759261
$ws.Range("A37").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B37").Value = "404号直流"
$ws.Range("C37").Value = 45967.556689814817
$ws.Range("A38").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B38").Value = "105号直流"
$ws.Range("C38").Value = 45967.577106481483
$ws.Range("A39").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B39").Value = "402号直流"
$ws.Range("C39").Value = 45967.577314814815
$ws.Range("A40").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B40").Value = "602号直流"
$ws.Range("C40").Value = 45967.592800925922
$ws.Range("A41").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B41").Value = "704号直流"
$ws.Range("C41").Value = 45967.599583333336
$ws.Range("A42").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B42").Value = "103号直流"
$ws.Range("C42").Value = 45967.609293981484
$ws.Range("A43").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B43").Value = "905号直流"
$ws.Range("C43").Value = 45967.624108796299
$ws.Range("A44").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B44").Value = "201号直流"
$ws.Range("C44").Value = 45967.63453703704
$ws.Range("A45").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B45").Value = "203号直流"
$ws.Range("C45").Value = 45967.646898148145
$ws.Range("A46").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B46").Value = "102号直流"
$ws.Range("C46").Value = 45967.664918981478
$ws.Range("A47").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B47").Value = "301号直流"
$ws.Range("C47").Value = 45967.678379629629
$ws.Range("A48").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B48").Value = "A03号直流"
$ws.Range("C48").Value = 45967.888356481482

# --- Rows 49-56 no longer have data in the refreshed export ---
$ws.Range("A49:D56").Value = ""

# --- Move the saved cursor/selection to match the refreshed view ---
$ws.Range("G10").Select()
